$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Client Id
$ws.Range("A2").Value = "FkJUh993"

# Update Candidate ID (numeric)
$ws.Range("B2").Value = 231010291

# Update User Name
$ws.Range("C2").Value = "qhonvob32"

# Update Exam Password
$ws.Range("D2").Value = "n%aN5&G8"

# Update First Name
$ws.Range("F2").Value = "kDjuRTMy"

# Update Last Name
$ws.Range("G2").Value = "Qfuf"
